# Weekly update: add two new "Fruta, Kiwi" price records for Feria
# Lagunitas de Puerto Montt at the top of the data block (row 73),
# pushing all the existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 73/74; Excel shifts rows 73:108 down to
# 75:110 (carrying their values/styles with them), which matches the
# target file exactly without needing to re-type the shifted rows.
$ws.Range("A73:A74").EntireRow.Insert()

# --- New row 73 : Kiwi Hayward, Especial ---
$ws.Cells.Item(73, 1).Value = 4
$ws.Cells.Item(73, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(73, 3).Value = "Los Lagos"
$ws.Cells.Item(73, 4).Value = 44460
$ws.Cells.Item(73, 5).Value = 10
$ws.Cells.Item(73, 6).Value = "Fruta"
$ws.Cells.Item(73, 7).Value = 100101
$ws.Cells.Item(73, 8).Value = "Berries"
$ws.Cells.Item(73, 9).Value = 100101007
$ws.Cells.Item(73, 10).Value = "Kiwi"
$ws.Cells.Item(73, 11).Value = "Hayward"
$ws.Cells.Item(73, 12).Value = "Especial"
$ws.Cells.Item(73, 13).Value = 200
$ws.Cells.Item(73, 14).Value = 20000
$ws.Cells.Item(73, 15).Value = 20000
$ws.Cells.Item(73, 16).Value = 20000
$ws.Cells.Item(73, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(73, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(73, 19).Value = 1333
$ws.Cells.Item(73, 20).Value = 15

# --- New row 74 : Kiwi Hayward, Primera ---
$ws.Cells.Item(74, 1).Value = 4
$ws.Cells.Item(74, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(74, 3).Value = "Los Lagos"
$ws.Cells.Item(74, 4).Value = 44460
$ws.Cells.Item(74, 5).Value = 10
$ws.Cells.Item(74, 6).Value = "Fruta"
$ws.Cells.Item(74, 7).Value = 100101
$ws.Cells.Item(74, 8).Value = "Berries"
$ws.Cells.Item(74, 9).Value = 100101007
$ws.Cells.Item(74, 10).Value = "Kiwi"
$ws.Cells.Item(74, 11).Value = "Hayward"
$ws.Cells.Item(74, 12).Value = "Primera"
$ws.Cells.Item(74, 13).Value = 100
$ws.Cells.Item(74, 14).Value = 14000
$ws.Cells.Item(74, 15).Value = 14000
$ws.Cells.Item(74, 16).Value = 14000
$ws.Cells.Item(74, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(74, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(74, 19).Value = 933
$ws.Cells.Item(74, 20).Value = 15
